# "Generate Report for Archive"
#
# The localization status report was regenerated: the outstanding
# zh-cn / de-de status changed from "Ready for handoff" to
# "In Translation", and the Status column on every sheet that shows it
# was narrowed to fit the new (shorter) text.
#
# "Ready for handoff" is shared by four cells across the three sheets
# (Overview!E2:F2, zh-cn!C2, de-de!C2) via the shared-strings table, so
# all four must be updated together for the single underlying string to
# change everywhere it is used.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("E2:F2").Value = "In Translation"
$zhcn.Range("C2").Value        = "In Translation"
$dede.Range("C2").Value        = "In Translation"

# Re-fit the Status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
